# Updated cryptos list on Thu Dec 28 22:57:54 UTC 2023 with GitHub Actions
# Refreshes price/volume figures in the coin table and fixes the ordering
# of three coin-name/link/price/volume rows that had gotten swapped
# (Uniswap/WrappedBTC, Kaspa/RenderToken, ARBITRUM/Algorand).
#
# Note: several "Price" column values look numeric (e.g. 325.36) but must
# stay stored as text, matching the sheet's existing inlineStr/text cells.
# A leading apostrophe is used (just like typing '325.36 into Excel) to
# force text entry for those values so Excel doesn't auto-convert them to
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.744.61'
$ws.Range('E2').Value = '  -1.77%  '

$ws.Range('D3').Value = '2.357.07'
$ws.Range('E3').Value = '  -0.45%  '

$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').Value = '''325.36'
$ws.Range('E5').Value = '  +2.46%  '

$ws.Range('D6').Value = '''103.86'
$ws.Range('E6').Value = '  -3.37%  '

$ws.Range('D7').Value = '''0.639'
$ws.Range('E7').Value = '  +0.15%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '''0.627'
$ws.Range('E9').Value = '  -1.32%  '

$ws.Range('D10').Value = '''40.43'
$ws.Range('E10').Value = '  -4.81%  '

$ws.Range('D11').Value = '''0.0923'
$ws.Range('E11').Value = '  -1.18%  '

$ws.Range('D12').Value = '''8.48'
$ws.Range('E12').Value = '  -2.34%  '

$ws.Range('D13').Value = '''1.00'
$ws.Range('E13').Value = '  -3.70%  '

$ws.Range('E14').Value = '  +0.15%  '

$ws.Range('D15').Value = '''16.20'
$ws.Range('E15').Value = '  -2.63%  '

$ws.Range('D16').Value = '2.708.23'
$ws.Range('E16').Value = '  -0.51%  '

$ws.Range('D17').Value = '2.345.49'
$ws.Range('E17').Value = '  -1.06%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '''7.98'
$ws.Range('E18').Value = '  +10.39%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '42.700.73'
$ws.Range('E19').Value = '  -1.85%  '

$ws.Range('E20').Value = '  -1.71%  '

$ws.Range('D21').Value = '''76.58'
$ws.Range('E21').Value = '  +1.77%  '

$ws.Range('D22').Value = '''3.70'
$ws.Range('E22').Value = '  +6.59%  '

$ws.Range('D23').Value = '''263.68'
$ws.Range('E23').Value = '  +1.94%  '

$ws.Range('D24').Value = '''2.32'
$ws.Range('E24').Value = '  -8.21%  '

$ws.Range('D25').Value = '''10.14'
$ws.Range('E25').Value = '  +8.91%  '

$ws.Range('E26').Value = '  +0.25%  '

$ws.Range('D27').Value = '''11.51'
$ws.Range('E27').Value = '  -3.81%  '

$ws.Range('D28').Value = '''22.88'
$ws.Range('E28').Value = '  +0.65%  '

$ws.Range('E29').Value = '  -1.32%  '

$ws.Range('D30').Value = '''175.50'
$ws.Range('E30').Value = '  +1.10%  '

$ws.Range('D31').Value = '''3.11'
$ws.Range('E31').Value = '  -3.32%  '

$ws.Range('D32').Value = '''0.0900'
$ws.Range('E32').Value = '  -2.04%  '

$ws.Range('D33').Value = '''35.47'
$ws.Range('E33').Value = '  -8.32%  '

$ws.Range('D34').Value = '''6.14'
$ws.Range('E34').Value = '  +3.30%  '

$ws.Range('D35').Value = '''0.133'
$ws.Range('E35').Value = '  +1.05%  '

$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '''4.57'
$ws.Range('E36').Value = '  -7.44%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '''0.110'
$ws.Range('E37').Value = '  +5.41%  '

$ws.Range('D38').Value = '''0.0358'
$ws.Range('E38').Value = '  -3.61%  '

$ws.Range('D39').Value = '''3.80'
$ws.Range('E39').Value = '  -6.36%  '

$ws.Range('D40').Value = '''2.81'
$ws.Range('E40').Value = '  +0.99%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '''0.238'
$ws.Range('E41').Value = '  +2.51%  '

$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''1.49'
$ws.Range('E42').Value = '  -0.65%  '

$ws.Range('D43').Value = '''70.12'
$ws.Range('E43').Value = '  -1.90%  '

$ws.Range('D44').Value = '''122.59'
$ws.Range('E44').Value = '  +9.69%  '

$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('D46').Value = '''92.86'
$ws.Range('E46').Value = '  +22.65%  '

$ws.Range('D47').Value = '''11.95'
$ws.Range('E47').Value = '  -4.92%  '

$ws.Range('D48').Value = '''5.54'
$ws.Range('E48').Value = '  -1.13%  '

$ws.Range('D49').Value = '''9.18'
$ws.Range('E49').Value = '  -1.23%  '

$ws.Range('E50').Value = '  -3.11%  '

$ws.Range('E51').Value = '  +0.35%  '
